# Applies the NORTH_DAKOTA_2018.xlsx cleanup edit:
#  - rename header columns to snake_case machine-friendly names
#  - normalize capitalization of a handful of municipality names
#    (lowercase connector words "de"/"y"/"el"/"la" -> capitalized)
#  - remove the trailing metadata/footer rows (74-78), shrinking the
#    used range from A1:D78 down to A1:D72

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize municipality name capitalization
$ws.Range("B22").Value = "Cuautepec De Hinojosa"
$ws.Range("B23").Value = "Tenango De Doria"
$ws.Range("B30").Value = "San Miguel El Alto"
$ws.Range("B49").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B50").Value = "Tlacolula De Matamoros"
$ws.Range("B57").Value = "San Pedro De La Cueva"

# 3. Remove the footer/metadata rows 74-78 (sample size, source, author, etc.)
$ws.Range("A74:A78").EntireRow.Delete()
